# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (header in G1 is "K"). Update rows 2-13
# with the newly regenerated K values.
$kValues = @{
    2  = 2
    3  = 2
    4  = 2
    5  = 1
    6  = 1
    7  = 3
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 3
    13 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
